$wb = $excel.ActiveWorkbook

# --- Swap the "Delete" / "Clone" entries on the Actions sheet (A6 <-> A7) ---
$actionsWs = $wb.Worksheets.Item("Actions")
$a6 = $actionsWs.Range("A6").Value2
$a7 = $actionsWs.Range("A7").Value2
$actionsWs.Range("A6").Value2 = $a7
$actionsWs.Range("A7").Value2 = $a6

# --- Move the selection away from the ExpenseRequest sheet's tab ---
$expenseWs = $wb.Worksheets.Item("ExpenseRequest")
$expenseWs.Range("E12").Select()

# --- Make the Actions sheet the active / selected sheet with A12 selected ---
$actionsWs.Activate()
$actionsWs.Range("A12").Select()
